# Adds 12 new match-result rows (10-21) to the game data sheet,
# reflecting newly implemented offensive/defensive heuristics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Cells.Item(10, 1).Value = "agent_random"
$ws.Cells.Item(10, 2).Value = "agent_minimax"
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = "agent_random"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "17.1794"
$ws.Cells.Item(10, 5).ClearFormats()
$ws.Cells.Item(10, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(10, 7).Value = "2 1 1 0 0 0`n1 1 2 1 0 0`n2 2 2 0 0 0`n1 1 2 2 1 0`n2 1 2 2 0 0`n1 2 2 1 1 0`n1 2 0 0 0 0"
$ws.Cells.Item(10, 8).Value = 5
$ws.Rows.Item(10).AutoFit()

# Row 11
$ws.Cells.Item(11, 1).Value = "agent_minimax"
$ws.Cells.Item(11, 2).Value = "agent_random"
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(11, 4).Value = "agent_minimax"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "8.4680"
$ws.Cells.Item(11, 5).ClearFormats()
$ws.Cells.Item(11, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(11, 7).Value = "2 1 2 1 2 0`n0 0 0 0 0 0`n2 1 1 1 1 0`n1 2 1 2 0 0`n2 0 0 0 0 0`n2 1 2 1 2 0`n1 0 0 0 0 0"
$ws.Cells.Item(11, 8).Value = 5
$ws.Rows.Item(11).AutoFit()

# Row 12
$ws.Cells.Item(12, 1).Value = "agent_minimax"
$ws.Cells.Item(12, 2).Value = "agent_minimax"
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(12, 4).Value = "agent_minimax"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "9.1750"
$ws.Cells.Item(12, 5).ClearFormats()
$ws.Cells.Item(12, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(12, 7).Value = "1 2 0 0 0 0`n1 2 0 0 0 0`n2 2 0 0 0 0`n1 2 1 0 0 0`n1 0 0 0 0 0`n0 0 0 0 0 0`n0 0 0 0 0 0"
$ws.Cells.Item(12, 8).Value = 5
$ws.Rows.Item(12).AutoFit()

# Row 13
$ws.Cells.Item(13, 1).Value = "agent_minimax"
$ws.Cells.Item(13, 2).Value = "agent_minimax"
$ws.Cells.Item(13, 3).Value = 2
$ws.Cells.Item(13, 4).Value = "agent_minimax"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "7.4280"
$ws.Cells.Item(13, 5).ClearFormats()
$ws.Cells.Item(13, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(13, 7).Value = "2 0 0 0 0 0`n1 2 2 0 0 0`n1 1 2 0 0 0`n1 1 1 2 0 0`n2 0 0 0 0 0`n0 0 0 0 0 0`n0 0 0 0 0 0"
$ws.Cells.Item(13, 8).Value = 5
$ws.Rows.Item(13).AutoFit()

# Row 14
$ws.Cells.Item(14, 1).Value = "agent_minimax"
$ws.Cells.Item(14, 2).Value = "agent_minimax"
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = "agent_minimax"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "16.1865"
$ws.Cells.Item(14, 5).ClearFormats()
$ws.Cells.Item(14, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(14, 7).Value = "2 1 0 0 0 0`n1 2 0 0 0 0`n2 1 1 1 1 0`n1 2 2 0 0 0`n2 0 0 0 0 0`n1 2 0 0 0 0`n2 1 2 1 0 0"
$ws.Cells.Item(14, 8).Value = 5
$ws.Rows.Item(14).AutoFit()

# Row 15
$ws.Cells.Item(15, 1).Value = "agent_minimax"
$ws.Cells.Item(15, 2).Value = "agent_minimax"
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 4).Value = "agent_minimax"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "15.8501"
$ws.Cells.Item(15, 5).ClearFormats()
$ws.Cells.Item(15, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(15, 7).Value = "2 1 2 2 1 2`n1 1 1 2 1 2`n1 0 0 0 0 0`n1 2 2 1 2 1`n2 1 1 1 2 1`n2 1 0 0 0 0`n1 2 2 2 1 2"
$ws.Cells.Item(15, 8).Value = 5
$ws.Rows.Item(15).AutoFit()

# Row 16
$ws.Cells.Item(16, 1).Value = "agent_minimax"
$ws.Cells.Item(16, 2).Value = "agent_minimax"
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(16, 4).Value = "agent_minimax"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "12.5153"
$ws.Cells.Item(16, 5).ClearFormats()
$ws.Cells.Item(16, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(16, 7).Value = "2 2 1 2 0 0`n1 0 0 0 0 0`n2 1 1 1 2 0`n1 2 1 2 1 2`n1 1 2 0 0 0`n2 1 2 2 2 1`n1 0 0 0 0 0"
$ws.Cells.Item(16, 8).Value = 5
$ws.Rows.Item(16).AutoFit()

# Row 17
$ws.Cells.Item(17, 1).Value = "agent_user"
$ws.Cells.Item(17, 2).Value = "agent_minimax"
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = "agent_user"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "61.7931"
$ws.Cells.Item(17, 5).ClearFormats()
$ws.Cells.Item(17, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(17, 7).Value = "2 0 0 0 0 0`n2 1 2 0 0 0`n1 2 1 2 0 0`n1 2 1 1 2 0`n1 1 2 1 1 0`n2 0 0 0 0 0`n0 0 0 0 0 0"
$ws.Cells.Item(17, 8).Value = 5
$ws.Rows.Item(17).AutoFit()

# Row 18
$ws.Cells.Item(18, 1).Value = "agent_minimax"
$ws.Cells.Item(18, 2).Value = "agent_minimax"
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).Value = "agent_minimax"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "7.0454"
$ws.Cells.Item(18, 5).ClearFormats()
$ws.Cells.Item(18, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(18, 7).Value = "2 0 0 0 0 0`n2 2 1 0 0 0`n2 1 2 0 0 0`n1 1 1 2 0 0`n2 2 1 0 0 0`n1 1 2 1 0 0`n0 0 0 0 0 0"
$ws.Cells.Item(18, 8).Value = 5
$ws.Rows.Item(18).AutoFit()

# Row 19
$ws.Cells.Item(19, 1).Value = "agent_minimax"
$ws.Cells.Item(19, 2).Value = "agent_minimax"
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = "agent_minimax"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "3.5081"
$ws.Cells.Item(19, 5).ClearFormats()
$ws.Cells.Item(19, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(19, 7).Value = "0 0 0 0 0 0`n2 1 1 2 2 0`n1 0 0 0 0 0`n1 1 2 0 0 0`n2 1 0 0 0 0`n2 1 0 0 0 0`n2 1 0 0 0 0"
$ws.Cells.Item(19, 8).Value = 5
$ws.Rows.Item(19).AutoFit()

# Row 20
$ws.Cells.Item(20, 1).Value = "agent_minimax"
$ws.Cells.Item(20, 2).Value = "agent_minimax"
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = "agent_minimax"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "9.7790"
$ws.Cells.Item(20, 5).ClearFormats()
$ws.Cells.Item(20, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(20, 7).Value = "1 2 2 2 1 2`n1 1 2 0 0 0`n2 1 2 0 0 0`n1 1 2 0 0 0`n1 2 2 1 2 0`n2 2 1 2 0 0`n1 1 2 1 1 1"
$ws.Cells.Item(20, 8).Value = 5
$ws.Rows.Item(20).AutoFit()

# Row 21
$ws.Cells.Item(21, 1).Value = "agent_minimax"
$ws.Cells.Item(21, 2).Value = "agent_minimax"
$ws.Cells.Item(21, 3).Value = 2
$ws.Cells.Item(21, 4).Value = "agent_minimax"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "16.7147"
$ws.Cells.Item(21, 5).ClearFormats()
$ws.Cells.Item(21, 6).Value = "h_block_fork|h_center_control"
$ws.Cells.Item(21, 7).Value = "1 1 2 2 1 2`n2 1 2 0 0 0`n2 2 1 2 1 1`n1 2 1 2 2 2`n1 1 1 2 1 2`n2 0 0 0 0 0`n1 2 1 1 2 1"
$ws.Cells.Item(21, 8).Value = 5
$ws.Rows.Item(21).AutoFit()

$ws.Range("A1").Select()
